# Insert a new data row at row 453 on the active sheet, pushing the existing
# rows 453:535 down to 454:536 (a new weekly Pomelo price observation for
# "Vega Modelo de Temuco" was added to the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 453 (and everything below it) down by one row.
$ws.Rows("453:453").Insert()

# Populate the newly inserted row 453 with the new record. Most of the
# categorical columns repeat the same product/market metadata as the
# surrounding rows.
$ws.Range("A453").Value = 10
$ws.Range("B453").Value = "Vega Modelo de Temuco"
$ws.Range("C453").Value = "La Araucanía"
$ws.Range("D453").Value = 45258
$ws.Range("E453").Value = 9
$ws.Range("F453").Value = "Fruta"
$ws.Range("G453").Value = 100102
$ws.Range("H453").Value = "Cítricos"
$ws.Range("I453").Value = 100102006
$ws.Range("J453").Value = "Pomelo"
$ws.Range("K453").Value = "Start Ruby"
$ws.Range("L453").Value = "Primera"
$ws.Range("M453").Value = 125
$ws.Range("N453").Value = 14000
$ws.Range("O453").Value = 14000
$ws.Range("P453").Value = 14000
$ws.Range("Q453").Value = "$/caja 14 kilos granel"
$ws.Range("R453").Value = "Región de O'Higgins"
$ws.Range("S453").Value = 1000
$ws.Range("T453").Value = 14
